$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("facility")
$ws.Activate()

# Update the H column formulas: F/3 -> F/5 (H2 is a standalone formula,
# H3:H11 is a shared formula group)
$ws.Range("H2").Formula = "=F2/5"
$ws.Range("H3:H11").Formula = "=F3/5"

# Highlight rows where A = 1, 5, 6, 8, 9 (rows 3, 7, 8, 10, 11) in yellow
$ws.Range("A3").Interior.Color = 65535
$ws.Range("A7").Interior.Color = 65535
$ws.Range("A8").Interior.Color = 65535
$ws.Range("A10").Interior.Color = 65535
$ws.Range("A11").Interior.Color = 65535

# Update the selected cell shown in the sheet view
$ws.Range("H17").Select()
